$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.875.54'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.216.76'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.18'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.79'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.213.47'
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.31'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.59'
$ws.Range("E14").Value = '  +4.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.749.82'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.222.97'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.920.60'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.60'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.33'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.07'
$ws.Range("E21").Value = '  +1.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").Value = '  -2.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.71'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.47'
$ws.Range("E24").Value = '  -0.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.60'
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.90'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.10'
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.89'
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.62'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.103'
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.46'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("E35").Value = '  -3.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").Value = '  +0.52%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.79'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0736'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0396'
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.78'
$ws.Range("E40").Value = '  +4.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '405.65'
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.16'
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.836.14'
$ws.Range("E44").Value = '  -7.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.257'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.17'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '129.41'
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.13'
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.89'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("E51").Value = '  -0.10%  '
